$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Shapefile Info" to "ast_config"
$ws.Name = "ast_config"
